# Add "Test case 2" to the Guru99 live project test-case sheet.
#
# The original sheet has a header row (row 1) and a single test case
# row (row 2). This change inserts a second test case (row 3) plus a
# handful of blank, pre-formatted rows below it (rows 4-15), matching
# the author's habit of pre-formatting rows ahead of future test cases.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row right after the existing test case. Inserting (rather
# than just writing into a fresh row) makes Excel clone the formatting of
# the row above, so the new row automatically gets the same per-column
# styles (centered/top aligned id column, wrapped text columns, etc.)
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).RowHeight = 140

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Verify that cost of product in list page and details page are equal"
$ws.Range("C3").Value = "1. Goto http://live.demoguru99.com`n2. Cick one mobile menu`n3. In the list of all mobile, read the cost of sony xperia mobile. note this value`n4. Click on Sony Xperia mobile`n5. Read the Sony xoperia mobile from detail page`n6. Compare value in step 3 & 5"
$ws.Range("E3").Value = "Product value in list and details page should be equal(`$100)"

# Pre-format a dozen more empty rows below the new test case, again by
# inserting (cloning formatting) rather than typing into blank rows.
for ($i = 4; $i -le 15; $i++) {
    $ws.Rows.Item($i).Insert()
    $ws.Rows.Item($i).RowHeight = 19
}

# Header row gains wrap-text (so the header stays readable once the
# columns below started wrapping multi-line step/expected-result text).
$ws.Range("A1:E1").WrapText = $true
$ws.Rows.Item(1).RowHeight = 20

$null = $ws.Range("C3").Select()
